$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Time column (A2:A5) with the new measured results.
$ws.Range("A2").Value = 3.39
$ws.Range("A3").Value = 12.34
$ws.Range("A4").Value = 186.32
$ws.Range("A5").Value = 66.01000000000001
